$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.468507333333333
$ws.Range("H2").Value = 4.405521999999999
$ws.Range("I2").Value = 0.005118279455112885
$ws.Range("J2").Value = 0.005118279455112885
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.356166333333333
$ws.Range("N2").Value = 4.068499
$ws.Range("O2").Value = 0.4438852433350371
$ws.Range("P2").Value = 0.4438852433350371
$ws.Range("Q2").Value = 1.991540205719778
$ws.Range("R2").Value = 17.923861851478
$ws.Range("S2").Value = 0.002271928721389504
$ws.Range("T2").Value = 0.002271928721389504

$ws.Range("G3").Value = 1.468507333333333
$ws.Range("H3").Value = 4.405521999999999
$ws.Range("I3").Value = 0.005118279455112885
$ws.Range("J3").Value = 0.005118279455112885
$ws.Range("M3").Value = 0.9609030000000001
$ws.Range("O3").Value = 0.3145120561487422
$ws.Range("P3").Value = 0.3145120561487422
$ws.Range("Q3").Value = 1.411093102122
$ws.Range("R3").Value = 12.699837919098
$ws.Range("S3").Value = 0.001609760595371417
$ws.Range("T3").Value = 0.001609760595371417

$ws.Range("G4").Value = 1.468507333333333
$ws.Range("H4").Value = 4.405521999999999
$ws.Range("I4").Value = 0.005118279455112885
$ws.Range("J4").Value = 0.005118279455112885
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.525837
$ws.Range("N4").Value = 1.577511
$ws.Range("O4").Value = 0.1721111038981938
$ws.Range("P4").Value = 0.1721111038981938
$ws.Range("Q4").Value = 0.7721954906379999
$ws.Range("R4").Value = 6.949759415741999
$ws.Range("S4").Value = 0.0008809127270789247
$ws.Range("T4").Value = 0.0008809127270789247

$ws.Range("G5").Value = 1.468507333333333
$ws.Range("H5").Value = 4.405521999999999
$ws.Range("I5").Value = 0.005118279455112885
$ws.Range("J5").Value = 0.005118279455112885
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.212312
$ws.Range("N5").Value = 0.6369359999999999
$ws.Range("O5").Value = 0.06949159661802674
$ws.Range("P5").Value = 0.06949159661802674
$ws.Range("Q5").Value = 0.3117817289546666
$ws.Range("R5").Value = 2.806035560592
$ws.Range("S5").Value = 0.0003556774112730383
$ws.Range("T5").Value = 0.0003556774112730383

$ws.Range("I6").Value = 0.9046276674881553
$ws.Range("J6").Value = 0.9046276674881553
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.356166333333333
$ws.Range("N6").Value = 4.068499
$ws.Range("O6").Value = 0.4438852433350371
$ws.Range("P6").Value = 0.4438852433350371
$ws.Range("Q6").Value = 351.9937484479203
$ws.Range("R6").Value = 3167.943736031282
$ws.Range("S6").Value = 0.4015508723105869
$ws.Range("T6").Value = 0.4015508723105868

$ws.Range("I7").Value = 0.9046276674881553
$ws.Range("J7").Value = 0.9046276674881553
$ws.Range("M7").Value = 0.9609030000000001
$ws.Range("O7").Value = 0.3145120561487422
$ws.Range("P7").Value = 0.3145120561487422
$ws.Range("Q7").Value = 249.4029239271181
$ws.Range("S7").Value = 0.2845163077507404
$ws.Range("T7").Value = 0.2845163077507404

$ws.Range("I8").Value = 0.9046276674881553
$ws.Range("J8").Value = 0.9046276674881553
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.525837
$ws.Range("N8").Value = 1.577511
$ws.Range("O8").Value = 0.1721111038981938
$ws.Range("P8").Value = 0.1721111038981938
$ws.Range("Q8").Value = 136.481294479322
$ws.Range("R8").Value = 1228.331650313898
$ws.Range("S8").Value = 0.1556964664682347
$ws.Range("T8").Value = 0.1556964664682347

$ws.Range("I9").Value = 0.9046276674881553
$ws.Range("J9").Value = 0.9046276674881553
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.212312
$ws.Range("N9").Value = 0.6369359999999999
$ws.Range("O9").Value = 0.06949159661802674
$ws.Range("P9").Value = 0.06949159661802674
$ws.Range("Q9").Value = 55.10570118400533
$ws.Range("R9").Value = 495.951310656048
$ws.Range("S9").Value = 0.06286402095859331
$ws.Range("T9").Value = 0.06286402095859331

$ws.Range("G10").Value = 0.5890733333333333
$ws.Range("H10").Value = 1.76722
$ws.Range("I10").Value = 0.002053133730501083
$ws.Range("J10").Value = 0.002053133730501083
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.356166333333333
$ws.Range("N10").Value = 4.068499
$ws.Range("O10").Value = 0.4438852433350371
$ws.Range("P10").Value = 0.4438852433350371
$ws.Range("Q10").Value = 0.7988814225311112
$ws.Range("R10").Value = 7.18993280278
$ws.Range("S10").Value = 0.0009113557655628458
$ws.Range("T10").Value = 0.0009113557655628457

$ws.Range("G11").Value = 0.5890733333333333
$ws.Range("H11").Value = 1.76722
$ws.Range("I11").Value = 0.002053133730501083
$ws.Range("J11").Value = 0.002053133730501083
$ws.Range("M11").Value = 0.9609030000000001
$ws.Range("O11").Value = 0.3145120561487422
$ws.Range("P11").Value = 0.3145120561487422
$ws.Range("Q11").Value = 0.5660423332200001
$ws.Range("R11").Value = 5.09438099898
$ws.Range("S11").Value = 0.0006457353111282332
$ws.Range("T11").Value = 0.0006457353111282332

$ws.Range("G12").Value = 0.5890733333333333
$ws.Range("H12").Value = 1.76722
$ws.Range("I12").Value = 0.002053133730501083
$ws.Range("J12").Value = 0.002053133730501083
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.525837
$ws.Range("N12").Value = 1.577511
$ws.Range("O12").Value = 0.1721111038981938
$ws.Range("P12").Value = 0.1721111038981938
$ws.Range("Q12").Value = 0.30975655438
$ws.Range("R12").Value = 2.78780898942
$ws.Range("S12").Value = 0.0003533671128071582
$ws.Range("T12").Value = 0.0003533671128071582

$ws.Range("G13").Value = 0.5890733333333333
$ws.Range("H13").Value = 1.76722
$ws.Range("I13").Value = 0.002053133730501083
$ws.Range("J13").Value = 0.002053133730501083
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.212312
$ws.Range("N13").Value = 0.6369359999999999
$ws.Range("O13").Value = 0.06949159661802674
$ws.Range("P13").Value = 0.06949159661802674
$ws.Range("Q13").Value = 0.1250673375466667
$ws.Range("R13").Value = 1.12560603792
$ws.Range("S13").Value = 0.0001426755410028457
$ws.Range("T13").Value = 0.0001426755410028457

$ws.Range("G14").Value = 25.306101
$ws.Range("H14").Value = 75.91830299999999
$ws.Range("I14").Value = 0.0882009193262308
$ws.Range("J14").Value = 0.0882009193262308
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.356166333333333
$ws.Range("N14").Value = 4.068499
$ws.Range("O14").Value = 0.4438852433350371
$ws.Range("P14").Value = 0.4438852433350371
$ws.Range("Q14").Value = 34.319282204133
$ws.Range("R14").Value = 308.873539837197
$ws.Range("S14").Value = 0.03915108653749794
$ws.Range("T14").Value = 0.03915108653749793

$ws.Range("G15").Value = 25.306101
$ws.Range("H15").Value = 75.91830299999999
$ws.Range("I15").Value = 0.0882009193262308
$ws.Range("J15").Value = 0.0882009193262308
$ws.Range("M15").Value = 0.9609030000000001
$ws.Range("O15").Value = 0.3145120561487422
$ws.Range("P15").Value = 0.3145120561487422
$ws.Range("Q15").Value = 24.316708369203
$ws.Range("R15").Value = 218.850375322827
$ws.Range("S15").Value = 0.02774025249150218
$ws.Range("T15").Value = 0.02774025249150218

$ws.Range("G16").Value = 25.306101
$ws.Range("H16").Value = 75.91830299999999
$ws.Range("I16").Value = 0.0882009193262308
$ws.Range("J16").Value = 0.0882009193262308
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.525837
$ws.Range("N16").Value = 1.577511
$ws.Range("O16").Value = 0.1721111038981938
$ws.Range("P16").Value = 0.1721111038981938
$ws.Range("Q16").Value = 13.306884231537
$ws.Range("R16").Value = 119.761958083833
$ws.Range("S16").Value = 0.01518035759007312
$ws.Range("T16").Value = 0.01518035759007312

$ws.Range("G17").Value = 25.306101
$ws.Range("H17").Value = 75.91830299999999
$ws.Range("I17").Value = 0.0882009193262308
$ws.Range("J17").Value = 0.0882009193262308
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.212312
$ws.Range("N17").Value = 0.6369359999999999
$ws.Range("O17").Value = 0.06949159661802674
$ws.Range("P17").Value = 0.06949159661802674
$ws.Range("Q17").Value = 5.372788915511999
$ws.Range("R17").Value = 48.35510023960799
$ws.Range("S17").Value = 0.00612922270715755
$ws.Range("T17").Value = 0.00612922270715755
